$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ML_Features")

# Shift RiskLabel out, insert new columns, then set final order:
# K1 = CyberRiskScore, L1 = DeviceTrustScore, M1 = RiskLabel
$ws.Range("K1").Value = "CyberRiskScore"
$ws.Range("L1").Value = "DeviceTrustScore"
$ws.Range("M1").Value = "RiskLabel"
